$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 (Q7)
$ws.Range("B9").Value = 0.4513396938152655
$ws.Range("C9").Value = 0.4513396938152655
$ws.Range("D9").Value = 0.2671320870154189
$ws.Range("E9").Value = 0.5168482243516165
$ws.Range("F9").Value = 0.275879468903711
$ws.Range("G9").Value = 6

# Row 10 (Q8)
$ws.Range("B10").Value = 0.1308504580670433
$ws.Range("C10").Value = 0.1308504580670433
$ws.Range("D10").Value = 0.0229929799516238
$ws.Range("E10").Value = 0.1516343627006221
$ws.Range("F10").Value = 0.09384405342323576
$ws.Range("G10").Value = 3

# Row 11 (Q9)
$ws.Range("B11").Value = -0.09392443396517081
$ws.Range("C11").Value = 0.09392443396517081
$ws.Range("D11").Value = 0.008821799295677731
$ws.Range("E11").Value = 0.09392443396517081
$ws.Range("F11").ClearContents()
$ws.Range("G11").Value = 1
